# Update "想去人数" (F) and "最低票价" (G) figures across sheets, reflecting
# refreshed scrape counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 69
$ws.Range("F3").Value = 126
$ws.Range("F4").Value = 2056
$ws.Range("F6").Value = 593
$ws.Range("F8").Value = 2064
$ws.Range("F9").Value = 10570
$ws.Range("F11").Value = 155
$ws.Range("F14").Value = 411
$ws.Range("F15").Value = 7470
$ws.Range("F18").Value = 224

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = 120

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 69
$ws.Range("F3").Value = 126
$ws.Range("F4").Value = 2056
$ws.Range("F6").Value = 593
$ws.Range("F7").Value = 22
$ws.Range("G7").Value = 120
$ws.Range("F9").Value = 2064
$ws.Range("F12").Value = 10570
$ws.Range("F14").Value = 155
$ws.Range("F17").Value = 411
$ws.Range("F18").Value = 7470
$ws.Range("F21").Value = 224
